$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.175.01'
$ws.Range('E2').Value = '  -6.70%  '
$ws.Range('D3').Value = '1.671.74'
$ws.Range('E3').Value = '  -4.38%  '
$ws.Range('E4').Value = '  +0.45%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.75'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5071'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -12.66%  '
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2637'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06318'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.43%  '
$ws.Range('E10').Value = '  -6.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07391'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.44%  '
$ws.Range('D12').Value = '1.674.30'
$ws.Range('E12').Value = '  -4.10%  '
$ws.Range('E13').Value = '  -3.82%  '
$ws.Range('E14').Value = '  -5.13%  '
$ws.Range('D15').Value = '1.898.63'
$ws.Range('E15').Value = '  -4.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000008531'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.36%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.61'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -12.77%  '
$ws.Range('D18').Value = '26.225.78'
$ws.Range('E18').Value = '  -6.40%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.960'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.80%  '
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('E21').Value = '  -4.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '186.71'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -9.27%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.180'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.05%  '
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.11'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.626'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.95%  '
$ws.Range('E27').Value = '  -5.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.310'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05740'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.00%  '
$ws.Range('E31').Value = '  -4.32%  '
$ws.Range('E32').Value = '  -6.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.492'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.667'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.5979'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.02%  '
$ws.Range('E37').Value = '  -2.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.633'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('D39').Value = '1.091.56'
$ws.Range('E39').Value = '  -3.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01598'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.906'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8611'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.41%  '
$ws.Range('E43').Value = '  -0.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').Value = '1.819.45'
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000112'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.03'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.002'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.054'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4316'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05202'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.61%  '
